$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Table grid: the 4 data-column widths go from 971 -> 970 twips
#    (21.5pt -> ... ; Word stores Width in points, so 970/20 = 48.5pt)
# -----------------------------------------------------------------
$tbl1 = $d.Tables(1)
for ($i = 2; $i -le $tbl1.Columns.Count; $i++) {
  $tbl1.Columns($i).Width = 48.5
}

# -----------------------------------------------------------------
# 2) "...contr" + " means contribution to national" + " " + "poverty."
#    -> merge into a single run " means contribution to national poverty."
#    (identical resulting text, so do a 2-step replace to force the
#    run rebuild/merge instead of a same-text no-op)
# -----------------------------------------------------------------
$rngA = $d.Content
$foundA = $rngA.Find.Execute(" means contribution to national poverty.")
if ($foundA) {
  $rngA.Text = " means contribution to national poverty.###TMP###"
  $rngA2 = $d.Content
  $foundA2 = $rngA2.Find.Execute(" means contribution to national poverty.###TMP###")
  if ($foundA2) {
    $rngA2.Text = " means contribution to national poverty."
  }
}

# -----------------------------------------------------------------
# 3) After "Marital status" insert a new run " at the" with the same
#    rFonts(cs)/szCs(bi) formatting as the surrounding paragraph text.
# -----------------------------------------------------------------
$rngB = $d.Content
$foundB = $rngB.Find.Execute("Marital status")
if ($foundB) {
  $rngB.Collapse(0)
  $rngB.Select()
  $sel = $word.Selection
  $sel.InsertAfter(" at the")
  $sel.MoveEnd(1, 7)
  $sel.Font.NameBi = "Times New Roman"
  $sel.Font.SizeBi = 12
}

# -----------------------------------------------------------------
# 4) "Source: Author's calculations from the different waves of the
#    UNPS." -> "Source: Author's calculations from the UNPS."
#    There are 4 occurrences; the first 3 are simple text swaps.
#    (Use Range.Text, not Find's ReplaceWith, so the straight
#    apostrophe isn't smart-quoted.)
# -----------------------------------------------------------------
$oldSentence = "Source: Author's calculations from the different waves of the UNPS."
$newSentence = "Source: Author's calculations from the UNPS."

for ($i = 0; $i -lt 3; $i++) {
  $rngC = $d.Content
  $foundC = $rngC.Find.Execute($oldSentence)
  if ($foundC) {
    $rngC.Text = $newSentence
  }
}

# -----------------------------------------------------------------
# 5) Final (4th) occurrence: split into two runs around the moved
#    "_GoBack" bookmark, and clear the bookmark out of the trailing
#    empty paragraph that used to host it.
# -----------------------------------------------------------------
$lastRng = $null
$searchRng = $d.Content
while ($true) {
  $f = $searchRng.Find.Execute($oldSentence)
  if (-not $f) { break }
  $lastRng = $d.Range($searchRng.Start, $searchRng.End)
  $searchRng.Collapse(0)
}

if ($lastRng -ne $null) {
  $lastRng.Text = $newSentence

  $prefix = "Source: Author's calculations from "
  $splitOffset = $lastRng.Start + $prefix.Length
  $bmRng = $d.Range($splitOffset, $splitOffset)
  $d.Bookmarks.Add("_GoBack", $bmRng)
}

Write-Output "done"
